$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new columns before column D (shifts old D:K to F:M)
$ws.Range("D:E").EntireColumn.Insert()

# Copy cell formatting (number format/font/style) from the shifted old columns (now F:G)
# into the new D:E columns, restricted to the row ranges that actually contain data
# (avoids creating phantom blank cells in label-only rows like 5, 6, 37, 79).
$fmtRanges = @(
    @{ From = "F7:G35";   To = "D7:E35" },
    @{ From = "F38:G77";  To = "D38:E77" },
    @{ From = "F80:G102"; To = "D80:E102" }
)
foreach ($pair in $fmtRanges) {
    $ws.Range($pair.From).Copy()
    $ws.Range($pair.To).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# Populate the two new columns (D = period ending 2018-12-31, E = period ending 2018-09-30)
# with the newly reported quarterly figures.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 257200
$ws.Range("E8").Value = 242100
$ws.Range("D9").Value = 208900
$ws.Range("E9").Value = 186300
$ws.Range("D10").Value = 48300
$ws.Range("E10").Value = 55800
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 4800
$ws.Range("E14").Value = 100
$ws.Range("D15").Value = 19500
$ws.Range("E15").Value = 17500
$ws.Range("D17").Value = 257800
$ws.Range("E17").Value = 225700
$ws.Range("D18").Value = -600
$ws.Range("E18").Value = 16400
$ws.Range("D20").Value = 41600
$ws.Range("E20").Value = 2900
$ws.Range("D21").Value = 60500
$ws.Range("E21").Value = 36700
$ws.Range("D22").Value = 12100
$ws.Range("E22").Value = 10700
$ws.Range("D23").Value = 28900
$ws.Range("E23").Value = 8600
$ws.Range("D24").Value = -2400
$ws.Range("E24").Value = 2800
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 31300
$ws.Range("E26").Value = 5800
$ws.Range("D27").Value = 29100
$ws.Range("E27").Value = 5000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -41600
$ws.Range("E32").Value = -2900
$ws.Range("D33").Value = 29100
$ws.Range("E33").Value = 5000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 29100
$ws.Range("E35").Value = 5000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 10400
$ws.Range("E41").Value = 27200
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 149500
$ws.Range("E43").Value = 157200
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 48800
$ws.Range("E45").Value = 42600
$ws.Range("D46").Value = 208700
$ws.Range("E46").Value = 227000
$ws.Range("D47").Value = 38000
$ws.Range("E47").Value = 42200
$ws.Range("D48").Value = 345700
$ws.Range("E48").Value = 285800
$ws.Range("D49").Value = 458700
$ws.Range("E49").Value = 313400
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 58300
$ws.Range("E52").Value = 42300
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1109300
$ws.Range("E54").Value = 910600
$ws.Range("D57").Value = 194100
$ws.Range("E57").Value = 164300
$ws.Range("D58").Value = 39300
$ws.Range("E58").Value = 33400
$ws.Range("D59").Value = 6100
$ws.Range("E59").Value = 5700
$ws.Range("D60").Value = 239500
$ws.Range("E60").Value = 203400
$ws.Range("D61").Value = 633000
$ws.Range("E61").Value = 553000
$ws.Range("D62").Value = 36500
$ws.Range("E62").Value = 33000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 982100
$ws.Range("E66").Value = 816300
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -117900
$ws.Range("E72").Value = -147100
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 127200
$ws.Range("E76").Value = 94300
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 29100
$ws.Range("E81").Value = 5000
$ws.Range("D83").Value = 19500
$ws.Range("E83").Value = 17500
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 28900
$ws.Range("E89").Value = 40900
$ws.Range("D91").Value = -54800
$ws.Range("E91").Value = -3300
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -65100
$ws.Range("E94").Value = -20500
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 19400
$ws.Range("E100").Value = -9400
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -16800
$ws.Range("E102").Value = 11000

# Row 91 (Capital Expenditures) received restated historical figures in addition to the
# newly reported columns, so fix up the shifted F:J cells to the restated values.
$ws.Range("F91").Value = -8300
$ws.Range("G91").Value = -5800
$ws.Range("H91").Value = -4700
$ws.Range("I91").Value = -13000
$ws.Range("J91").Value = -6800

Write-Output "edit complete"
